$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.767.92"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "3.707.81"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'599.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'168.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.90%  "
$ws.Range("D7").Value = "3.706.35"
$ws.Range("E7").Value = "  -2.43%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("D11").Value = "'6.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("D13").Value = "'38.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "'0.0000244"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").Value = "4.326.74"
$ws.Range("E15").Value = "  -2.69%  "
$ws.Range("D16").Value = "3.703.75"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "68.660.04"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'17.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.63%  "
$ws.Range("D21").Value = "'494.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").Value = "'9.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.74%  "
$ws.Range("D24").Value = "'84.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "'12.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "'7.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("D33").Value = "'31.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.57%  "
$ws.Range("D34").Value = "3.845.14"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("D36").Value = "3.646.27"
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("E40").Value = "  -2.75%  "
$ws.Range("D41").Value = "'0.322"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'433.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").Value = "'49.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "'2.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "
$ws.Range("D46").Value = "'8.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'40.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").Value = "'142.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "2.735.15"
$ws.Range("E51").Value = "  -3.33%  "
